$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "586.79" -> 586.79) are written via a Text number format so they stay
# plain text, matching the source data (inline strings), then the format is
# reset to Normal so no stray style survives on the cell.

$ws.Range("D2").Value = '63.666.73'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '3.136.86'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.136.56'
$ws.Range("E8").Value = '  +0.99%  '
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("E10").Value = '  +5.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.73'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("E12").Value = '  -2.32%  '
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.17%  '
$ws.Range("D15").Value = '3.663.74'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("E16").Value = '  -1.75%  '
$ws.Range("D17").Value = '3.144.60'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '63.580.49'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '463.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("E26").Value = '  +2.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.94%  '
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.73%  '
$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").Value = '0.0₃0846'
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("E37").Value = '  -5.08%  '
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '440.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("D44").Value = '2.901.86'
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("E46").Value = '  -2.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '37.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.31%  '
